# Auto-generated edit script applying numeric updates from the commit diff.
# Updates static market-price / profit figures across 8 worksheets (Titan_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1227
$ws.Range("I32").Value = 1139.8
$ws.Range("J32").Value = 1336
$ws.Range("K32").Value = 1139.8
$ws.Range("L32").Value = 1336
$ws.Range("M32").Value = -813.8
$ws.Range("N32").Value = -1988

# Row 98
$ws.Range("H98").Value = 622462.9
$ws.Range("I98").Value = 931653.0600000001
$ws.Range("J98").Value = 4082.5
$ws.Range("K98").Value = 931653.0600000001
$ws.Range("L98").Value = 4082.5
$ws.Range("M98").Value = -930155.0600000001
$ws.Range("N98").Value = -7078.5

# Row 122
$ws.Range("H122").Value = 622462.9
$ws.Range("I122").Value = 931653.0600000001
$ws.Range("J122").Value = 4082.5
$ws.Range("K122").Value = 2794959.18
$ws.Range("L122").Value = 12247.5
$ws.Range("M122").Value = -2792509.18
$ws.Range("N122").Value = -17147.5

# Row 132
$ws.Range("H132").Value = 367046.22
$ws.Range("I132").Value = 434323.75
$ws.Range("J132").Value = 53084.332
$ws.Range("K132").Value = 1302971.25
$ws.Range("L132").Value = 159252.996
$ws.Range("M132").Value = -1300441.25
$ws.Range("N132").Value = -164312.996

# Row 137
$ws.Range("H137").Value = 1185.9412
$ws.Range("I137").Value = 654
$ws.Range("J137").Value = 1277.6552
$ws.Range("K137").Value = 1962
$ws.Range("L137").Value = 3832.9656
$ws.Range("M137").Value = 588
$ws.Range("N137").Value = -8932.9656

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2872.7256
$ws.Range("I32").Value = 2980.5588
$ws.Range("J32").Value = 2657.0588
$ws.Range("K32").Value = 2980.5588
$ws.Range("L32").Value = 2657.0588
$ws.Range("M32").Value = -2693.5588
$ws.Range("N32").Value = -3231.0588

# Row 61
$ws.Range("H61").Value = 2837.05
$ws.Range("I61").Value = 1894
$ws.Range("J61").Value = 5666.2
$ws.Range("K61").Value = 1894
$ws.Range("L61").Value = 5666.2
$ws.Range("M61").Value = -1682
$ws.Range("N61").Value = -6090.2

# Row 74
$ws.Range("H74").Value = 1094.5416
$ws.Range("I74").Value = 1120.6923
$ws.Range("J74").Value = 1063.6364
$ws.Range("K74").Value = 1120.6923
$ws.Range("L74").Value = 1063.6364
$ws.Range("M74").Value = -246.6922999999999
$ws.Range("N74").Value = -2811.6364

# Row 77
$ws.Range("H77").Value = 1094.5416
$ws.Range("I77").Value = 1120.6923
$ws.Range("J77").Value = 1063.6364
$ws.Range("K77").Value = 5603.461499999999
$ws.Range("L77").Value = 5318.182000000001
$ws.Range("M77").Value = -1235.461499999999
$ws.Range("N77").Value = -14054.182

# Row 122
$ws.Range("H122").Value = 3670.5
$ws.Range("I122").Value = 4425
$ws.Range("J122").Value = 2916
$ws.Range("K122").Value = 13275
$ws.Range("L122").Value = 8748
$ws.Range("M122").Value = -10825
$ws.Range("N122").Value = -13648

# Row 132
$ws.Range("H132").Value = 2610.5112
$ws.Range("I132").Value = 2061.0715
$ws.Range("J132").Value = 3515.4707
$ws.Range("K132").Value = 6183.2145
$ws.Range("L132").Value = 10546.4121
$ws.Range("M132").Value = -3653.2145
$ws.Range("N132").Value = -15606.4121

# Row 136
$ws.Range("H136").Value = 2837.05
$ws.Range("I136").Value = 1894
$ws.Range("J136").Value = 5666.2
$ws.Range("K136").Value = 5682
$ws.Range("L136").Value = 16998.6
$ws.Range("M136").Value = -3132
$ws.Range("N136").Value = -22098.6

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1900.2559
$ws.Range("I134").Value = 1296.6
$ws.Range("K134").Value = 3889.8
$ws.Range("M134").Value = -1354.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 325.16666
$ws.Range("I16").Value = 334.4
$ws.Range("J16").Value = 279
$ws.Range("K16").Value = 334.4
$ws.Range("L16").Value = 279
$ws.Range("M16").Value = -47.39999999999998
$ws.Range("N16").Value = -853

# Row 58
$ws.Range("H58").Value = 1828.7667
$ws.Range("I58").Value = 1329.5883
$ws.Range("J58").Value = 2481.5386
$ws.Range("K58").Value = 1329.5883
$ws.Range("L58").Value = 2481.5386
$ws.Range("M58").Value = -1126.5883
$ws.Range("N58").Value = -2887.5386

# Row 99
$ws.Range("H99").Value = 6946652
$ws.Range("J99").Value = 1542.5
$ws.Range("L99").Value = 1542.5
$ws.Range("N99").Value = -4538.5

# Row 107
$ws.Range("H107").Value = 709.4706
$ws.Range("I107").Value = 405.16666
$ws.Range("J107").Value = 875.4545000000001
$ws.Range("K107").Value = 405.16666
$ws.Range("L107").Value = 875.4545000000001
$ws.Range("M107").Value = 1514.83334
$ws.Range("N107").Value = -4715.4545

# Row 113
$ws.Range("H113").Value = 325.16666
$ws.Range("I113").Value = 334.4
$ws.Range("J113").Value = 279
$ws.Range("K113").Value = 334.4
$ws.Range("L113").Value = 279
$ws.Range("M113").Value = 1835.6
$ws.Range("N113").Value = -4619

# Row 122
$ws.Range("H122").Value = 1003.2143
$ws.Range("I122").Value = 810.3333
$ws.Range("J122").Value = 1147.875
$ws.Range("K122").Value = 2430.9999
$ws.Range("L122").Value = 3443.625
$ws.Range("M122").Value = 19.0001000000002
$ws.Range("N122").Value = -8343.625

# Row 126
$ws.Range("H126").Value = 6946652
$ws.Range("J126").Value = 1542.5
$ws.Range("L126").Value = 4627.5
$ws.Range("N126").Value = -9567.5

# Row 132
$ws.Range("H132").Value = 3438.32
$ws.Range("I132").Value = 1926.9231
$ws.Range("J132").Value = 5075.6665
$ws.Range("K132").Value = 5780.7693
$ws.Range("L132").Value = 15226.9995
$ws.Range("M132").Value = -3250.7693
$ws.Range("N132").Value = -20286.9995

# Row 136
$ws.Range("H136").Value = 1828.7667
$ws.Range("I136").Value = 1329.5883
$ws.Range("J136").Value = 2481.5386
$ws.Range("K136").Value = 3988.7649
$ws.Range("L136").Value = 7444.6158
$ws.Range("M136").Value = -1438.7649
$ws.Range("N136").Value = -12544.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1462.8182
$ws.Range("I5").Value = 1025
$ws.Range("J5").Value = 1560.1111
$ws.Range("K5").Value = 3075
$ws.Range("L5").Value = 4680.3333
$ws.Range("M5").Value = -2963
$ws.Range("N5").Value = -4904.3333

# Row 12
$ws.Range("H12").Value = 24.095238
$ws.Range("I12").Value = 9.111110999999999
$ws.Range("J12").Value = 35.333332
$ws.Range("K12").Value = 27.333333
$ws.Range("L12").Value = 105.999996
$ws.Range("M12").Value = 145.666667
$ws.Range("N12").Value = -451.999996

# Row 14
$ws.Range("H14").Value = 52.666668
$ws.Range("I14").Value = 52.666668
$ws.Range("K14").Value = 158.000004
$ws.Range("M14").Value = 14.99999600000001

# Row 23
$ws.Range("H23").Value = 1445.5555
$ws.Range("I23").Value = 2242.2
$ws.Range("J23").Value = 449.75
$ws.Range("K23").Value = 6726.599999999999
$ws.Range("L23").Value = 1349.25
$ws.Range("M23").Value = -6491.599999999999
$ws.Range("N23").Value = -1819.25

# Row 33
$ws.Range("H33").Value = 158
$ws.Range("I33").Value = 110.28571
$ws.Range("J33").Value = 325
$ws.Range("K33").Value = 661.71426
$ws.Range("L33").Value = 1950
$ws.Range("M33").Value = -378.71426
$ws.Range("N33").Value = -2516

# Row 113
$ws.Range("H113").Value = 400.9355
$ws.Range("I113").Value = 367.44446
$ws.Range("J113").Value = 414.63635
$ws.Range("K113").Value = 1102.33338
$ws.Range("L113").Value = 1243.90905
$ws.Range("M113").Value = 1067.66662
$ws.Range("N113").Value = -5583.90905

# Row 131
$ws.Range("H131").Value = 937.1900000000001
$ws.Range("I131").Value = 515
$ws.Range("J131").Value = 964.1383
$ws.Range("K131").Value = 1545
$ws.Range("L131").Value = 2892.4149
$ws.Range("M131").Value = 3495
$ws.Range("N131").Value = -12972.4149

# Row 132
$ws.Range("H132").Value = 825.86365
$ws.Range("I132").Value = 577.6667
$ws.Range("J132").Value = 997.6923
$ws.Range("K132").Value = 5199.0003
$ws.Range("L132").Value = 8979.2307
$ws.Range("M132").Value = -2669.0003
$ws.Range("N132").Value = -14039.2307

# Row 135
$ws.Range("H135").Value = 1462.8182
$ws.Range("I135").Value = 1025
$ws.Range("J135").Value = 1560.1111
$ws.Range("K135").Value = 9225
$ws.Range("L135").Value = 14040.9999
$ws.Range("M135").Value = -6690
$ws.Range("N135").Value = -19110.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2224202.2
$ws.Range("I122").Value = 11111111
$ws.Range("J122").Value = 2475
$ws.Range("K122").Value = 33333333
$ws.Range("L122").Value = 7425
$ws.Range("M122").Value = -33330883
$ws.Range("N122").Value = -12325

# Row 126
$ws.Range("H126").Value = 2148.5
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2310.625
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 6931.875
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -11871.875

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 682.6667
$ws.Range("I16").Value = 629.2
$ws.Range("K16").Value = 629.2
$ws.Range("M16").Value = -459.2

# Row 61
$ws.Range("H61").Value = 7196.6665
$ws.Range("I61").Value = 9900
$ws.Range("J61").Value = 1790
$ws.Range("K61").Value = 9900
$ws.Range("L61").Value = 1790
$ws.Range("M61").Value = -9698
$ws.Range("N61").Value = -2194

# Row 113
$ws.Range("H113").Value = 7196.6665
$ws.Range("I113").Value = 9900
$ws.Range("J113").Value = 1790
$ws.Range("K113").Value = 9900
$ws.Range("L113").Value = 1790
$ws.Range("M113").Value = -7730
$ws.Range("N113").Value = -6130

# Row 122
$ws.Range("H122").Value = 3498.3333
$ws.Range("J122").Value = 3498.3333
$ws.Range("L122").Value = 10494.9999
$ws.Range("N122").Value = -15394.9999

# Row 132
$ws.Range("H132").Value = 4874.643
$ws.Range("I132").Value = 3874.75
$ws.Range("J132").Value = 5274.6
$ws.Range("K132").Value = 11624.25
$ws.Range("L132").Value = 15823.8
$ws.Range("M132").Value = -9094.25
$ws.Range("N132").Value = -20883.8

# Row 136
$ws.Range("H136").Value = 6321.136
$ws.Range("I136").Value = 2845.3
$ws.Range("J136").Value = 9217.666999999999
$ws.Range("K136").Value = 8535.900000000001
$ws.Range("L136").Value = 27653.001
$ws.Range("M136").Value = -5985.900000000001
$ws.Range("N136").Value = -32753.001

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 84508.336
$ws.Range("I126").Value = 125776.875
$ws.Range("J126").Value = 1971.25
$ws.Range("K126").Value = 377330.625
$ws.Range("L126").Value = 5913.75
$ws.Range("M126").Value = -374860.625
$ws.Range("N126").Value = -10853.75

# Row 132
$ws.Range("H132").Value = 11365205
$ws.Range("I132").Value = 14706850
$ws.Range("J132").Value = 3614.5
$ws.Range("K132").Value = 44120550
$ws.Range("L132").Value = 10843.5
$ws.Range("M132").Value = -44118020
$ws.Range("N132").Value = -15903.5
